$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1550.8298
$ws.Range("J17").Value = 1582
$ws.Range("L17").Value = 4746
$ws.Range("N17").Value = -5082

$ws.Range("H98").Value = 2095.7856
$ws.Range("I98").Value = 1076.4546
$ws.Range("K98").Value = 1076.4546
$ws.Range("M98").Value = 421.5454

$ws.Range("H122").Value = 2095.7856
$ws.Range("I122").Value = 1076.4546
$ws.Range("K122").Value = 3229.3638
$ws.Range("M122").Value = -779.3638000000001

$ws.Range("H132").Value = 3907.9
$ws.Range("I132").Value = 2763.25
$ws.Range("K132").Value = 8289.75
$ws.Range("M132").Value = -5759.75

$ws.Range("H136").Value = 49999.832
$ws.Range("J136").Value = 49999.832
$ws.Range("L136").Value = 49999.832
$ws.Range("N136").Value = -60199.832

$ws.Range("H137").Value = 3617.6
$ws.Range("J137").Value = 6150
$ws.Range("L137").Value = 18450
$ws.Range("N137").Value = -23550

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3444.1177
$ws.Range("I88").Value = 3845.6
$ws.Range("J88").Value = 2870.5715
$ws.Range("K88").Value = 3845.6
$ws.Range("L88").Value = 2870.5715
$ws.Range("M88").Value = -3439.6
$ws.Range("N88").Value = -3682.5715

$ws.Range("H91").Value = 3444.1177
$ws.Range("I91").Value = 3845.6
$ws.Range("J91").Value = 2870.5715
$ws.Range("K91").Value = 3845.6
$ws.Range("L91").Value = 2870.5715
$ws.Range("M91").Value = -2441.6
$ws.Range("N91").Value = -5678.5715

$ws.Range("H132").Value = 4326.768
$ws.Range("I132").Value = 1443.1471
$ws.Range("K132").Value = 4329.4413
$ws.Range("M132").Value = -1799.4413

$ws.Range("H137").Value = 58333.168
$ws.Range("J137").Value = 58333.168
$ws.Range("L137").Value = 58333.168
$ws.Range("N137").Value = -68533.16800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 4092
$ws.Range("J54").Value = 4092
$ws.Range("L54").Value = 4092
$ws.Range("N54").Value = -5060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 732.1
$ws.Range("I7").Value = 734.0769
$ws.Range("J7").Value = 728.4286
$ws.Range("K7").Value = 734.0769
$ws.Range("L7").Value = 728.4286
$ws.Range("M7").Value = -621.0769
$ws.Range("N7").Value = -954.4286

$ws.Range("H16").Value = 3598.8333
$ws.Range("I16").Value = 3518.8
$ws.Range("J16").Value = 3999
$ws.Range("K16").Value = 3518.8
$ws.Range("L16").Value = 3999
$ws.Range("M16").Value = -3231.8
$ws.Range("N16").Value = -4573

$ws.Range("H62").Value = 3624.5
$ws.Range("I62").Value = 3624.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3624.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3000.5
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3624.5
$ws.Range("I65").Value = 3624.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 18122.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -15002.5
$ws.Range("N65").ClearContents()

$ws.Range("H107").Value = 478.42856
$ws.Range("I107").Value = 483.16666
$ws.Range("K107").Value = 483.16666
$ws.Range("M107").Value = 1436.83334

$ws.Range("H113").Value = 3598.8333
$ws.Range("I113").Value = 3518.8
$ws.Range("J113").Value = 3999
$ws.Range("K113").Value = 3518.8
$ws.Range("L113").Value = 3999
$ws.Range("M113").Value = -1348.8
$ws.Range("N113").Value = -8339

$ws.Range("H132").Value = 4178.5
$ws.Range("I132").Value = 3160.8096
$ws.Range("J132").Value = 6553.1113
$ws.Range("K132").Value = 9482.4288
$ws.Range("L132").Value = 19659.3339
$ws.Range("M132").Value = -6952.4288
$ws.Range("N132").Value = -24719.3339

$ws.Range("H134").Value = 3824.2327
$ws.Range("I134").Value = 2881.08
$ws.Range("K134").Value = 8643.24
$ws.Range("M134").Value = -6108.24

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 2500
$ws.Range("I28").Value = 2500
$ws.Range("K28").Value = 7500
$ws.Range("M28").Value = -7268

$ws.Range("H97").Value = 1085.2
$ws.Range("J97").Value = 1030.75
$ws.Range("L97").Value = 3092.25
$ws.Range("N97").Value = -4084.25

$ws.Range("H132").Value = 4725
$ws.Range("I132").Value = 500
$ws.Range("J132").Value = 5570
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 50130
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -55190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1003925.9
$ws.Range("I80").Value = 1004500.8
$ws.Range("J80").Value = 1003351
$ws.Range("K80").Value = 1004500.8
$ws.Range("L80").Value = 1003351
$ws.Range("M80").Value = -1003502.8
$ws.Range("N80").Value = -1005347

$ws.Range("H83").Value = 1003925.9
$ws.Range("I83").Value = 1004500.8
$ws.Range("J83").Value = 1003351
$ws.Range("K83").Value = 5022504
$ws.Range("L83").Value = 5016755
$ws.Range("M83").Value = -5017512
$ws.Range("N83").Value = -5026739

$ws.Range("H113").Value = 1024682.9
$ws.Range("I113").Value = 1275990.1
$ws.Range("J113").Value = 19454
$ws.Range("K113").Value = 1275990.1
$ws.Range("L113").Value = 19454
$ws.Range("M113").Value = -1273820.1
$ws.Range("N113").Value = -23794

$ws.Range("H126").Value = 3184.8667
$ws.Range("I126").Value = 2201
$ws.Range("J126").Value = 4309.2856
$ws.Range("K126").Value = 6603
$ws.Range("L126").Value = 12927.8568
$ws.Range("M126").Value = -4133
$ws.Range("N126").Value = -17867.8568

$ws.Range("H132").Value = 403345.47
$ws.Range("J132").Value = 4904.1113
$ws.Range("L132").Value = 14712.3339
$ws.Range("N132").Value = -19772.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1053.1818
$ws.Range("I55").Value = 433.7143
$ws.Range("K55").Value = 433.7143
$ws.Range("M55").Value = -260.7143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6301.3335
$ws.Range("I62").Value = 4001
$ws.Range("K62").Value = 4001
$ws.Range("M62").Value = -3377

$ws.Range("H65").Value = 6301.3335
$ws.Range("I65").Value = 4001
$ws.Range("K65").Value = 20005
$ws.Range("M65").Value = -16885

$ws.Range("H122").Value = 27782254
$ws.Range("I122").Value = 43481810
$ws.Range("K122").Value = 130445430
$ws.Range("M122").Value = -130442980

$ws.Range("H126").Value = 3845
$ws.Range("I126").Value = 3741.6667
$ws.Range("K126").Value = 11225.0001
$ws.Range("M126").Value = -8755.000100000001
